$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the prefab path strings to reflect the new Object Pool (tag -> ID) scheme
$ws.Range("B4").Value = "Projectile/Prefabs/Orange Explosion"
$ws.Range("B5").Value = "Projectile/Prefabs/Poison"

# Move the active selection to B6
$ws.Range("B6").Select()
